$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '27.722.11'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  +5.93%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.734.02'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  +4.63%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '227.71'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  +3.80%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.5457'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  +3.64%  '
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  -0.11%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2742'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  +2.17%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06721'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  +5.31%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '21.92'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  +6.19%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07787'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  +1.22%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '4.697'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  +1.72%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.972.60'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  +4.65%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '1.692.53'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  +3.81%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.5985'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  +6.01%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.0₅8426'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  +1.94%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '69.21'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  +5.19%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '27.727.63'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  +6.04%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '226.23'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  +18.46%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '4.830'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  +2.97%  '
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  -0.14%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '10.89'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  +4.97%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '6.218'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  -0.11%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '147.68'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  -1.28%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.738'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  +13.77%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.1252'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  +3.97%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '7.475'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  +2.54%  '
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  +6.44%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.05717'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  +1.18%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.312'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  +2.63%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.705'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  +5.93%  '
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  +3.97%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.681'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  +6.11%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.9784'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  +2.83%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.853'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  +1.80%  '
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  +1.25%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.6006'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  +3.84%  '
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  +3.96%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '5.935'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  -0.82%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.8511'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  +2.05%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.050.54'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  +2.14%  '
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  -0.10%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '101.79'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  +0.22%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.879.16'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  +4.68%  '
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  +9.34%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '59.69'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  +1.85%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '8.312'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  +3.11%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.4426'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  +1.95%  '
$ws.Range("B50").NumberFormat = "@"
$ws.Range("B50").Value = 'Frax'
$ws.Range("C50").NumberFormat = "@"
$ws.Range("C50").Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.005'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  -0.01%  '
$ws.Range("B51").NumberFormat = "@"
$ws.Range("B51").Value = 'Cronos'
$ws.Range("C51").NumberFormat = "@"
$ws.Range("C51").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.05339'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  -0.35%  '
